$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "289.98"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-4.26%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "30.80"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-4.15%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.942"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.58%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07130"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-9.43%"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-13.35%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.682"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-1.98%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-2.80%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8974"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-3.12%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1648"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-6.35%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07600"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-2.39%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08102"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-5.86%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03049"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-3.61%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1003"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.21%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001502"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.78%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005722"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.64%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.470"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.09%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.082"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.70%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3277"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.02%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1296"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.59%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.037"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-5.88%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.63%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.38%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.91%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-10.27%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001249"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.03%"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-7.09%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04367"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-9.12%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007328"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.24%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1305"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-4.47%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002006"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-15.00%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009117"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-12.80%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005991"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-2.40%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.04%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.247"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "173.92%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-3.19%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.04%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0001999"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.04%"
